# Updated data in light of fix_baddata_wip
#
# The sheet "master_baddata" used the sentinel value -999 (numeric) and the
# shared string "CONFLICT!" to flag rows/cells with bad or conflicting data.
# This pass replaces those placeholders with the corrected/real values that
# were recovered, using either real numbers or the text markers "BAD" /
# "Maybe" / "M" as appropriate.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master_baddata")

# Row 2
$ws.Range("N2").Value = "BAD"

# Row 6
$ws.Range("C6").Value = "BAD"
$ws.Range("E6").Value = "BAD"
$ws.Range("F6").Value = "BAD"
$ws.Range("G6").Value = "BAD"
$ws.Range("I6").Value = "BAD"
$ws.Range("J6").Value = "BAD"
$ws.Range("K6").Value = "BAD"
$ws.Range("N6").Value = "BAD"

# Row 7
$ws.Range("N7").Value = 9

# Rows 8-23: gender column (AV) recovered to "M"
$ws.Range("AV8").Value  = "M"
$ws.Range("AV9").Value  = "M"
$ws.Range("AV10").Value = "M"
$ws.Range("AV11").Value = "M"
$ws.Range("AV12").Value = "M"
$ws.Range("AV13").Value = "M"
$ws.Range("AV14").Value = "M"
$ws.Range("AV15").Value = "M"
$ws.Range("AV16").Value = "M"
$ws.Range("AV17").Value = "M"
$ws.Range("AV18").Value = "M"
$ws.Range("AV19").Value = "M"
$ws.Range("AV20").Value = "M"
$ws.Range("AV21").Value = "M"
$ws.Range("AV22").Value = "M"
$ws.Range("AV23").Value = "M"

# Row 24
$ws.Range("V24").Value  = 4
$ws.Range("W24").Value  = 3
$ws.Range("X24").Value  = 5
$ws.Range("Y24").Value  = 3
$ws.Range("Z24").Value  = 5
$ws.Range("AD24").Value = 4
$ws.Range("AF24").Value = 3
$ws.Range("AH24").Value = 5
$ws.Range("AI24").Value = 3
$ws.Range("AK24").Value = 5
$ws.Range("AV24").Value = "M"

# Rows 25-37: gender column
$ws.Range("AV25").Value = "M"
$ws.Range("AV26").Value = "M"
$ws.Range("AV27").Value = "M"
$ws.Range("AV28").Value = "M"
$ws.Range("AV29").Value = "M"
$ws.Range("AV30").Value = "M"
$ws.Range("AV31").Value = "M"
$ws.Range("AV32").Value = "M"
$ws.Range("AV33").Value = "M"
$ws.Range("AV34").Value = "M"
$ws.Range("AV35").Value = "M"
$ws.Range("AV36").Value = "M"
$ws.Range("AV37").Value = "M"

# Row 38
$ws.Range("AG38").Value = 1
$ws.Range("AV38").Value = "M"

# Row 39
$ws.Range("AR39").Value = 18
$ws.Range("AS39").Value = 5
$ws.Range("AT39").Value = 1
$ws.Range("AV39").Value = "M"

# Rows 40-49: gender column
$ws.Range("AV40").Value = "M"
$ws.Range("AV41").Value = "M"
$ws.Range("AV42").Value = "M"
$ws.Range("AV43").Value = "M"
$ws.Range("AV44").Value = "M"
$ws.Range("AV45").Value = "M"
$ws.Range("AV46").Value = "M"
$ws.Range("AV47").Value = "M"
$ws.Range("AV48").Value = "M"
$ws.Range("AV49").Value = "M"

# Row 50
$ws.Range("W50").Value  = "BAD"
$ws.Range("AV50").Value = "M"

# Rows 51-57: gender column (plus S54)
$ws.Range("AV51").Value = "M"
$ws.Range("AV52").Value = "M"
$ws.Range("AV53").Value = "M"
$ws.Range("S54").Value  = 3
$ws.Range("AV54").Value = "M"
$ws.Range("AV55").Value = "M"
$ws.Range("AV56").Value = "M"
$ws.Range("AV57").Value = "M"

# Row 58
$ws.Range("B58").Value  = "BAD"
$ws.Range("C58").Value  = "BAD"
$ws.Range("D58").Value  = "BAD"
$ws.Range("E58").Value  = "BAD"
$ws.Range("F58").Value  = "BAD"
$ws.Range("G58").Value  = "BAD"
$ws.Range("H58").Value  = "BAD"
$ws.Range("O58").Value  = 3
$ws.Range("P58").Value  = 1
$ws.Range("Q58").Value  = 6
$ws.Range("R58").Value  = 3
$ws.Range("V58").Value  = 6
$ws.Range("W58").Value  = 6
$ws.Range("X58").Value  = 2
$ws.Range("Y58").Value  = 3
$ws.Range("AA58").Value = 4
$ws.Range("AB58").Value = 5
$ws.Range("AC58").Value = 5
$ws.Range("AD58").Value = 2
$ws.Range("AE58").Value = 2
$ws.Range("AF58").Value = 3
$ws.Range("AG58").Value = 4
$ws.Range("AI58").Value = 3
$ws.Range("AJ58").Value = 5
$ws.Range("AK58").Value = 5
$ws.Range("AL58").Value = 4
$ws.Range("AM58").Value = 6
$ws.Range("AN58").Value = 6
$ws.Range("AO58").Value = 5
$ws.Range("AP58").Value = 3
$ws.Range("AR58").Value = 19
$ws.Range("AU58").Value = 3.82

# Row 59
$ws.Range("F59").Value  = 1
$ws.Range("I59").Value  = 0
$ws.Range("K59").Value  = 2
$ws.Range("L59").Value  = "Maybe"
$ws.Range("M59").Value  = "Maybe"
$ws.Range("R59").Value  = "BAD"
$ws.Range("T59").Value  = "BAD"
$ws.Range("U59").Value  = "BAD"
$ws.Range("V59").Value  = "BAD"
$ws.Range("W59").Value  = "BAD"
$ws.Range("Y59").Value  = "BAD"
$ws.Range("Z59").Value  = "BAD"
$ws.Range("AA59").Value = "BAD"
$ws.Range("AB59").Value = "BAD"
$ws.Range("AC59").Value = "BAD"
$ws.Range("AD59").Value = "BAD"
$ws.Range("AE59").Value = "BAD"
$ws.Range("AF59").Value = "BAD"
$ws.Range("AG59").Value = "BAD"
$ws.Range("AH59").Value = "BAD"
$ws.Range("AI59").Value = "BAD"
$ws.Range("AJ59").Value = "BAD"
$ws.Range("AK59").Value = "BAD"
$ws.Range("AL59").Value = "BAD"
$ws.Range("AM59").Value = "BAD"
$ws.Range("AO59").Value = "BAD"
$ws.Range("AP59").Value = "BAD"
$ws.Range("AQ59").Value = "BAD"
$ws.Range("AR59").Value = "BAD"
$ws.Range("AS59").Value = "BAD"
$ws.Range("AT59").Value = "BAD"
$ws.Range("AU59").Value = "BAD"
